{"js": "// The document's final paragraph is a stray caption line that describes\n// Python/SPARQLWrapper code which is no longer part of the document (the\n// author \"did a new render\" and dropped this trailing paragraph). Remove\n// that whole paragraph, leaving the preceding image paragraph as the last\n// paragraph before the section break.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = \"The below Python code uses SPARQLWrapper to retrieve data from Wikidata based on a SPARQL query.\";\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === target) {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document's final paragraph is a stray caption line that describes\n# Python/SPARQLWrapper code which is no longer part of the document (the\n# author \"did a new render\" and dropped this trailing paragraph). Remove\n# that whole paragraph, leaving the preceding image paragraph as the last\n# paragraph before the section break.\n$d = $word.ActiveDocument\n$target = \"The below Python code uses SPARQLWrapper to retrieve data from Wikidata based on a SPARQL query.\"\n\n# Walk backwards so deleting a paragraph doesn't disturb the index of the\n# ones we still need to visit.\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($text -eq $target) {\n        $p.Range.Delete()\n    }\n}\n"}
